$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!B2,B3,C2,C3 ; zh-cn!C2,C3 ; de-de!C2,C3 all share this text)
# ---------------------------------------------------------------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.UsedRange.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# ---------------------------------------------------------------------------
# 2) Populate "Latest Target File" (F) and "Latest Handback File" (G) columns
#    for the zh-cn and de-de detail sheets, with matching hyperlinks, and
#    stamp the "Latest Handback DateTime" (H) column.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn, row 2 (2d98ff5e-...)
$wsZh.Range("F2").Value = "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/ee75db98b6586b3e82827ede45811b9f8b3e34f0/e2e/2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.md", "", "", "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.md") | Out-Null
$wsZh.Range("G2").Value = "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.c0f2915d94b55e7a4ac10a39e43729f360aad0dc.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca639dd9e54440cdee1bb644d9b873142683e700/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.c0f2915d94b55e7a4ac10a39e43729f360aad0dc.zh-cn.xlf", "", "", "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.c0f2915d94b55e7a4ac10a39e43729f360aad0dc.zh-cn.xlf") | Out-Null
$wsZh.Range("H2").Value = "2016-03-13 04:12:54"

# zh-cn, row 3 (43da8692-...)
$wsZh.Range("F3").Value = "43da8692-b4e1-43c0-8a20-d2b084996d03.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/ee75db98b6586b3e82827ede45811b9f8b3e34f0/e2e/43da8692-b4e1-43c0-8a20-d2b084996d03.md", "", "", "43da8692-b4e1-43c0-8a20-d2b084996d03.md") | Out-Null
$wsZh.Range("G3").Value = "43da8692-b4e1-43c0-8a20-d2b084996d03.ff1d4235747e5a79cd8c35cab65aee8c20797f87.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca639dd9e54440cdee1bb644d9b873142683e700/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/43da8692-b4e1-43c0-8a20-d2b084996d03.ff1d4235747e5a79cd8c35cab65aee8c20797f87.zh-cn.xlf", "", "", "43da8692-b4e1-43c0-8a20-d2b084996d03.ff1d4235747e5a79cd8c35cab65aee8c20797f87.zh-cn.xlf") | Out-Null
$wsZh.Range("H3").Value = "2016-03-13 04:12:54"

# de-de, row 2 (2d98ff5e-...)
$wsDe.Range("F2").Value = "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/ee75db98b6586b3e82827ede45811b9f8b3e34f0/e2e/2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.md", "", "", "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.md") | Out-Null
$wsDe.Range("G2").Value = "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.c0f2915d94b55e7a4ac10a39e43729f360aad0dc.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66e7d9953cfa5042feb3b01454fb61812f1ca68f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.c0f2915d94b55e7a4ac10a39e43729f360aad0dc.de-de.xlf", "", "", "2d98ff5e-d7d2-4c99-956e-d2e0599f7bc4.c0f2915d94b55e7a4ac10a39e43729f360aad0dc.de-de.xlf") | Out-Null
$wsDe.Range("H2").Value = "2016-03-13 04:13:00"

# de-de, row 3 (43da8692-...)
$wsDe.Range("F3").Value = "43da8692-b4e1-43c0-8a20-d2b084996d03.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/ee75db98b6586b3e82827ede45811b9f8b3e34f0/e2e/43da8692-b4e1-43c0-8a20-d2b084996d03.md", "", "", "43da8692-b4e1-43c0-8a20-d2b084996d03.md") | Out-Null
$wsDe.Range("G3").Value = "43da8692-b4e1-43c0-8a20-d2b084996d03.ff1d4235747e5a79cd8c35cab65aee8c20797f87.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66e7d9953cfa5042feb3b01454fb61812f1ca68f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/43da8692-b4e1-43c0-8a20-d2b084996d03.ff1d4235747e5a79cd8c35cab65aee8c20797f87.de-de.xlf", "", "", "43da8692-b4e1-43c0-8a20-d2b084996d03.ff1d4235747e5a79cd8c35cab65aee8c20797f87.de-de.xlf") | Out-Null
$wsDe.Range("H3").Value = "2016-03-13 04:13:00"
